$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.315233826637268
$ws.Range("B1").Value = 1.979480624198914
$ws.Range("C1").Value = 2.754829168319702
$ws.Range("D1").Value = 3.786533117294312
$ws.Range("E1").Value = 1.023110032081604
